$d = $word.ActiveDocument

# wdReplaceAll = 2
$wdReplaceAll = 2

# --- Title paragraph: "Answers:" " " "Definite" " " "integration" -> "Answers: Definite integration"
$titlePara = $d.Paragraphs(1).Range
$titlePara.Find.ClearFormatting()
$titlePara.Find.Execute("Answers: Definite integration", $true, $false, $false, $false, $false, $true, 1, $false, `
    "Answers: Definite integration", $wdReplaceAll)

# --- Author paragraph: "Donald" " " "Campbell" -> "Donald Campbell"
$authorPara = $d.Paragraphs(2).Range
$authorPara.Find.ClearFormatting()
$authorPara.Find.Execute("Donald Campbell", $true, $false, $false, $false, $false, $true, 1, $false, `
    "Donald Campbell", $wdReplaceAll)

# --- Abstract paragraph: "Answers" " " "to" " " "questions" ... -> single run
$abstractPara = $d.Paragraphs(4).Range
$abstractPara.Find.ClearFormatting()
$abstractPara.Find.Execute("Answers to questions relating to the guide on definite integration.", $true, $false, $false, $false, $false, $true, 1, $false, `
    "Answers to questions relating to the guide on definite integration.", $wdReplaceAll)
